# Benchmark template: add the missing "idle_iterations" configuration row
# (Description / Name / Value / Flag) as a new row appended to Table2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 ("iterations", an Integer-flagged row) already carries the exact
# Description/Name/Value/Flag formatting (fonts, wrap, fill + border on the
# Flag cell) that the new row needs, so clone its formatting down to row 13
# before filling in the new content.
$ws.Range("A4:D4").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)  # xlPasteFormats

# New entry: "Number of measurements for taking idle current" / idle_iterations / 20
$ws.Range("A13").Value = "Number of measurements for taking idle current"
$ws.Range("B13").Value = "idle_iterations"
$ws.Range("C13").Value = 20

# Leave the selection where the author ended up after entering the row.
$ws.Range("C15").Select()
